$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, centered, bordered) from an existing
# header cell onto the two new header cells before setting their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new I0 / IF columns (rows 2-11)
$i0 = @(9, 9, 5, 8, 9, 4, 8, 4, 5, 5)
$if = @(9, 9, 6, 9, 9, 4, 9, 4, 5, 5)

for ($r = 0; $r -lt 10; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $i0[$r]
    $ws.Cells.Item($row, 10).Value = $if[$r]
}
